$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D3";   Value = -8.6944 },
    @{ Cell = "E3";   Value = 16.01950000000001 },
    @{ Cell = "C12";  Value = -11.3788 },
    @{ Cell = "D14";  Value = -7.609999999999999 },
    @{ Cell = "D26";  Value = -9.063799999999993 },
    @{ Cell = "E30";  Value = 15.87839999999999 },
    @{ Cell = "D31";  Value = -9.052999999999988 },
    @{ Cell = "C32";  Value = -13.3833 },
    @{ Cell = "D35";  Value = -8.389299999999995 },
    @{ Cell = "C36";  Value = -12.47920000000001 },
    @{ Cell = "D37";  Value = -7.942999999999996 },
    @{ Cell = "C38";  Value = -12.74349999999999 },
    @{ Cell = "E44";  Value = 16.66489999999999 },
    @{ Cell = "D45";  Value = -7.682500000000002 },
    @{ Cell = "C46";  Value = -14.57049999999999 },
    @{ Cell = "C54";  Value = -13.1804 },
    @{ Cell = "C55";  Value = -13.7249 },
    @{ Cell = "D57";  Value = -8.240499999999999 },
    @{ Cell = "E58";  Value = 16.20420000000001 },
    @{ Cell = "C67";  Value = -11.0105 },
    @{ Cell = "C69";  Value = -12.1628 },
    @{ Cell = "C72";  Value = -11.5109 },
    @{ Cell = "E84";  Value = 16.51839999999999 },
    @{ Cell = "E89";  Value = 17.20180000000002 },
    @{ Cell = "C91";  Value = -10.61659999999999 },
    @{ Cell = "E91";  Value = 17.93700000000002 },
    @{ Cell = "E92";  Value = 18.05980000000002 },
    @{ Cell = "C99";  Value = -13.41599999999999 },
    @{ Cell = "D100"; Value = -8.112800000000004 },
    @{ Cell = "D102"; Value = -7.731100000000001 },
    @{ Cell = "E102"; Value = 16.4599 }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
